$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged) - update values
$ws.Range("B3").Value = 0.993279785718383
$ws.Range("C3").Value = 0.9934722483273073
$ws.Range("D3").Value = 0.9869030749314215

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9861716140379135
$ws.Range("C4").Value = 0.9866269962310491
$ws.Range("D4").Value = 0.9638027419938622

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9961388843886264
$ws.Range("C5").Value = 0.9961659672376908
$ws.Range("D5").Value = 0.9954775103081985
